$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2 and R2 to the nearest integer
$ws.Range("Q2").Value = 625419
$ws.Range("R2").Value = 6542623

# Clear the time cells Z2 (Starttid) and AB2 (Sluttid) entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
